# Updated code quality rules and rel rating
#
# The "BannedPaths" rule row is removed from its old position (row 35) and
# re-added (renamed to "BannedPath", singular) at the end of the
# "Bug"/"Blocker" block (row 40) with an updated (Critical) severity and no
# tags. Rows 36-39 each shift up by one to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 35: was "BannedPaths" -> becomes what used to be row 36
$ws.Range("A35").Value = "CloudServiceIncompatibleWorkflowProcess"
$ws.Range("B35").Value = "Usage of Cloud Service Incompatible Workflow Processes"
$ws.Range("C35").Value = "Bug"
$ws.Range("D35").Value = "Blocker"
$ws.Range("E35").Value = "aem,cloud-service-compatibility"

# Row 36: becomes what used to be row 37
$ws.Range("A36").Value = "IndexType"
$ws.Range("B36").Value = "Custom Search Index Definition Nodes Must Use the Index Type lucene"
$ws.Range("C36").Value = "Bug"
$ws.Range("D36").Value = "Blocker"
$ws.Range("E36").Value = "aem,cloud-service-compatibility"

# Row 37: becomes what used to be row 38
$ws.Range("A37").Value = "IndexAsyncProperty"
$ws.Range("B37").Value = "Custom Lucene Oak Indexes must not be synchronous"
$ws.Range("C37").Value = "Bug"
$ws.Range("D37").Value = "Blocker"
$ws.Range("E37").Value = "aem,cloud-service-compatibility"

# Row 38: becomes what used to be row 39
$ws.Range("A38").Value = "IndexTikaNode"
$ws.Range("B38").Value = "Custom Oak indexes must have a tika configuration"
$ws.Range("C38").Value = "Bug"
$ws.Range("D38").Value = "Blocker"
$ws.Range("E38").Value = "aem,cloud-service-compatibility"

# Row 39: becomes what used to be row 40
$ws.Range("A39").Value = "IndexDamAssetLucene"
$ws.Range("B39").Value = "Index customizations of the damAssetLucene Oak index should be properly structured."
$ws.Range("C39").Value = "Bug"
$ws.Range("D39").Value = "Blocker"
$ws.Range("E39").Value = "aem,cloud-service-compatibility"

# Row 40: re-added "BannedPath" (renamed, singular) with updated severity and no tags
$ws.Range("A40").Value = "BannedPath"
$ws.Range("B40").Value = "Customer packages should not install content under /libs"
$ws.Range("C40").Value = "Bug"
$ws.Range("D40").Value = "Critical"
$ws.Range("E40").ClearContents()

# Update the saved selection to match the author's final cursor position
$ws.Range("A37").Select()
